$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Trees and hedges information" section (originally a single row, 151)
# gains a second field ("Trees on adjacent land"), so a new row is inserted
# immediately below the existing row 151. This pushes every row from the old
# 152 down to 153, etc. (old 152-162 -> new 153-163).
# ---------------------------------------------------------------------------
$ws.Rows("152:152").Insert()

# Row 151 keeps its A/B header text ("Trees and hedges information" / "Details
# of trees and/or hedges...") but now also carries the first of the two field
# rows ("Trees on site").
$ws.Range("C151").Value = "Trees on site"
$ws.Range("G151").Value = "Whether trees or hedges are present on the proposed development site"
$ws.Range("H151").Value = "boolean"
$ws.Range("I151").Value = "MUST"

# The newly inserted row 152 holds the second field ("Trees on adjacent
# land"); A152/B152 stay blank since the header is merged down from row 151.
$ws.Range("C152").Value = "Trees on adjacent land"
$ws.Range("G152").Value = "Whether trees or hedges on land adjacent to the proposed development site could influence the development or might be important as part of the local landscape character"
$ws.Range("H152").Value = "boolean"
$ws.Range("I152").Value = "MUST"

# The A/B header cells for "Trees and hedges information" now span both
# rows 151:152 (previously just row 151).
$ws.Range("A151:A152").Merge()
$ws.Range("B151:B152").Merge()
